$d = $word.ActiveDocument

# "Allow setting custom DPI through plugin": the [[imageWithDPI]] merge field
# gains a plugin-style modifier so authors can pin an explicit DPI, turning
# [[imageWithDPI]] into [[imageWithDPI]:dpi(160)]
$d.Content.Find.Execute("[[imageWithDPI]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[[imageWithDPI]:dpi(160)]", 2)
